$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A20").Value = "NCT00780494"
$ws.Range("C20").Value = "hybrid"

$ws.Range("A21").Value = "NCT01474382"
$ws.Range("C21").Value = "closed"

$ws.Range("A19").Value = "NCT02440789"
$ws.Range("C19").Value = "bronze"

$ws.Range("A22").Value = "NCT02494024"
$ws.Range("C22").Value = "hybrid"

$ws.Range("C4").Value = "bronze"

$ws.Range("C25").Select()

$wb.Save()
